$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'332.09"
$ws.Range("E2").Value = "'0.56%"
$ws.Range("D3").Value = "'45.37"
$ws.Range("E3").Value = "'2.33%"
$ws.Range("D4").Value = "'5.553"
$ws.Range("E4").Value = "'1.11%"
$ws.Range("D5").Value = "'0.08375"
$ws.Range("E5").Value = "'4.36%"
$ws.Range("E6").Value = "'0.19%"
$ws.Range("D7").Value = "'0.9894"
$ws.Range("E7").Value = "'3.71%"
$ws.Range("D9").Value = "'0.1184"
$ws.Range("E9").Value = "'3.30%"
$ws.Range("D10").Value = "'0.1933"
$ws.Range("E10").Value = "'0.99%"
$ws.Range("D11").Value = "'9.560"
$ws.Range("E11").Value = "'-5.86%"
$ws.Range("D12").Value = "'0.09968"
$ws.Range("E12").Value = "'0.57%"
$ws.Range("D13").Value = "'0.04666"
$ws.Range("E13").Value = "'-3.55%"
$ws.Range("E14").Value = "'-0.50%"
$ws.Range("D15").Value = "'0.001295"
$ws.Range("E15").Value = "'1.59%"
$ws.Range("D16").Value = "'0.005918"
$ws.Range("E16").Value = "'-0.65%"
$ws.Range("D17").Value = "'3.393"
$ws.Range("E17").Value = "'0.59%"
$ws.Range("D18").Value = "'4.429"
$ws.Range("E18").Value = "'0.57%"
$ws.Range("D20").Value = "'0.1353"
$ws.Range("E20").Value = "'-1.92%"
$ws.Range("D22").Value = "'0.04132"
$ws.Range("E22").Value = "'1.27%"
$ws.Range("D23").Value = "'0.001292"
$ws.Range("E23").Value = "'1.51%"
$ws.Range("D24").Value = "'0.004531"
$ws.Range("E24").Value = "'5.53%"
$ws.Range("D25").Value = "'0.0001301"
$ws.Range("E25").Value = "'8.46%"
$ws.Range("E26").Value = "'0.01%"
$ws.Range("D38").Value = "'0.02704"
$ws.Range("E38").Value = "'4.42%"
$ws.Range("D39").Value = "'0.05772"
$ws.Range("E39").Value = "'-0.35%"
$ws.Range("D40").Value = "'0.007899"
$ws.Range("E40").Value = "'4.47%"
$ws.Range("D41").Value = "'0.1432"
$ws.Range("E41").Value = "'2.04%"
$ws.Range("D42").Value = "'0.007933"
$ws.Range("E42").Value = "'8.36%"
$ws.Range("D43").Value = "'0.002022"
$ws.Range("E43").Value = "'0.37%"
$ws.Range("D44").Value = "'0.008951"
$ws.Range("E44").Value = "'-1.24%"
$ws.Range("D45").Value = "'0.3412"
$ws.Range("D46").Value = "'0.00007059"
$ws.Range("E46").Value = "'0.75%"
$ws.Range("E48").Value = "'0.28%"
$ws.Range("D49").Value = "'0.003535"
$ws.Range("E49").Value = "'0.13%"
$ws.Range("D50").Value = "'0.003380"
$ws.Range("E50").Value = "'-4.43%"
$ws.Range("D51").Value = "'0.00002103"
